$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 173-175 with revised figures ---

# Row 173 (01-04-2021)
$ws.Range("B173").Value = 16995
$ws.Range("C173").Value = 21058
$ws.Range("D173").Value = 65216
$ws.Range("E173").Value = 44159
$ws.Range("J173").Value = -6634
$ws.Range("K173").Value = 10911
$ws.Range("L173").Value = 17546

# Row 174 (01-05-2021)
$ws.Range("B174").Value = 12148
$ws.Range("C174").Value = 19256
$ws.Range("D174").Value = 67343
$ws.Range("E174").Value = 48087
$ws.Range("J174").Value = -5770
$ws.Range("K174").Value = 10988
$ws.Range("L174").Value = 16759
$ws.Range("X174").Value = 3397

# Row 175 (01-06-2021)
$ws.Range("B175").Value = 13283
$ws.Range("C175").Value = 20143
$ws.Range("D175").Value = 63581
$ws.Range("E175").Value = 43438
$ws.Range("J175").Value = -5290
$ws.Range("K175").Value = 11214
$ws.Range("L175").Value = 16504

# --- Add new row 176 (01-07-2021) ---
# Force the date-like text to be stored as a plain string (matching column A's
# existing entries) instead of letting Excel auto-convert it to a date serial.
$ws.Range("A176").NumberFormat = "@"
$ws.Range("A176").Value = "01-07-2021"
$ws.Range("A176").Style = "Normal"
$ws.Range("B176").Value = 13180
$ws.Range("C176").Value = 20505
$ws.Range("D176").Value = 61029
$ws.Range("E176").Value = 40524
$ws.Range("F176").Value = 56351
$ws.Range("G176").Value = 15767
$ws.Range("H176").Value = 35065
$ws.Range("I176").Value = 19298
$ws.Range("J176").Value = -4407
$ws.Range("K176").Value = 11324
$ws.Range("L176").Value = 15731
$ws.Range("M176").Value = 8301
$ws.Range("N176").Value = 8630
$ws.Range("O176").Value = 329
$ws.Range("P176").Value = 164
$ws.Range("Q176").Value = 1320
$ws.Range("R176").Value = 1156
$ws.Range("S176").Value = 815
$ws.Range("T176").Value = 2126
$ws.Range("U176").Value = 1312
$ws.Range("V176").Value = -135
$ws.Range("W176").Value = 2563
$ws.Range("X176").Value = 2698
$ws.Range("Y176").Value = -7325
$ws.Range("Z176").Value = 56758
$ws.Range("AA176").Value = 64083
